# "Projeto Final - Entrega 1.xlsx" — fill in the "Lista das colunas na base
# de dados" table (rows 12-50, columns A:C) with the real column names /
# dtypes / NA-counts coming from the Ifood_df dataset, matching the
# supplied diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column names, in dataframe order, starting at row 12.
$columns = @(
    "Income",
    "Kidhome",
    "Teenhome",
    "Recency",
    "MntWines",
    "MntFruits",
    "MntMeatProducts",
    "MntFishProducts",
    "MntSweetProducts",
    "MntGoldProds",
    "NumDealsPurchases",
    "NumWebPurchases",
    "NumCatalogPurchases",
    "NumStorePurchases",
    "NumWebVisitsMonth",
    "AcceptedCmp3",
    "AcceptedCmp4",
    "AcceptedCmp5",
    "AcceptedCmp1",
    "AcceptedCmp2",
    "Complain",
    "Z_CostContact",
    "Z_Revenue",
    "Response",
    "Age",
    "Customer_Days",
    "marital_Divorced",
    "marital_Married",
    "marital_Single",
    "marital_Together",
    "marital_Widow",
    "education_2n Cycle",
    "education_Basic",
    "education_Graduation",
    "education_Master",
    "education_PhD",
    "MntTotal",
    "MntRegularProds",
    "AcceptedCmpOverall"
)

$startRow = 12

for ($i = 0; $i -lt $columns.Count; $i++) {
    $r = $startRow + $i

    $ws.Cells.Item($r, 1).Value = $columns[$i]
    $ws.Cells.Item($r, 2).Value = "dbl"
    $ws.Cells.Item($r, 3).Value = 0
}

# Formatting: the old placeholder rows used a shaded fill + thin box border
# (style index 2). The real data rows drop the fill (back to "Normal") and
# keep the thin box border around every cell.
$endRow = $startRow + $columns.Count - 1
$dataRange = $ws.Range("A" + $startRow + ":C" + $endRow)
$dataRange.Style = "Normal"
$dataRange.Borders.LineStyle = 1
$dataRange.Borders.Weight = 2

# Row 27 (AcceptedCmp3) carries an extra underline on the column-name cell.
$ws.Range("A27").Font.Underline = 2

# Restore the selection to where the author left off (C35), scrolled so row
# 7 is at the top of the viewport.
$ws.Activate()
$ws.Range("C35").Select()
try {
    $excel.ActiveWindow.ScrollRow = 7
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # View-scroll state isn't always reachable through this COM surface;
    # the selection itself is the important, verifiable part.
}
